$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.384.21'
$ws.Range('E2').Value = '  -0.68%  '

$ws.Range('D3').Value = '1.859.90'
$ws.Range('E3').Value = '  -1.26%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.64%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.21%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4737'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.38%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2743'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.82%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06436'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.22%  '

$ws.Range('D10').Value = '1.906.13'
$ws.Range('E10').Value = '  +0.60%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07459'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.23%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.34'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.38%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.997'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.83%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '85.62'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.66%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6340'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.39%  '

$ws.Range('D16').Value = '30.351.26'
$ws.Range('E16').Value = '  -0.55%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9996'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.05%  '

$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.78'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.61%  '

$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '230.18'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.52%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007428'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.76%  '

$ws.Range('D21').Value = '2.099.14'
$ws.Range('E21').Value = '  -5.14%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.18%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.009'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.17%  '

$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.997'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.06%  '

$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.265'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.45%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.35'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.66%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.97'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.64%  '

$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.895'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.65%  '

$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1045'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +7.16%  '

$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.401'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.01%  '

$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.155'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.16%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.935'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.84%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04935'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.16%  '

$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.165'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.90%  '

$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7263'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.43%  '

$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.15%  '

$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.701'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.13%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01880'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.83%  '

$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.652'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.20%  '

$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9172'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.68%  '

$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.972'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.87%  '

$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '106.10'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.36%  '

$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.000'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.09%  '

$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4118'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.69%  '

$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.580'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.35%  '

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.117'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.40%  '

$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '61.07'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.78%  '

$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1202'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.69%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.675'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.64%  '

$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.54'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.12%  '

$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.408'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.32%  '

Write-Host "done"